$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.745.18"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").Value = "2.220.03"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'241.31"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'74.74"
$ws.Range("E7").Value = "  +4.49%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  +4.18%  "
$ws.Range("D10").Value = "'41.46"
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "'54.71"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "'6.92"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "2.549.21"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'14.72"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "2.216.23"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "42.647.27"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "'70.82"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'5.94"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").Value = "'9.78"
$ws.Range("E23").Value = "  -7.38%  "
$ws.Range("D24").Value = "'229.73"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("E25").Value = "  +8.76%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("D28").Value = "'3.38"
$ws.Range("E28").Value = "  -6.78%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'172.98"
$ws.Range("E30").Value = "  +3.90%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'36.90"
$ws.Range("E31").Value = "  +20.16%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.09"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "'20.28"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Value = "'0.0799"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'5.30"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("E39").Value = "  +9.24%  "
$ws.Range("D40").Value = "'12.83"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'60.90"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").Value = "'8.64"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").Value = "'0.0992"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'99.58"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'2.30"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.436"
$ws.Range("E50").Value = "  +20.90%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.14"
$ws.Range("E51").Value = "  -1.30%  "
